# The workbook "excelFile_WithValidationError.xlsx" originally held sample
# employee rows whose Email column (and A/B/Dipartimento columns) had been
# scrambled/offset by one record versus the Nome/Cognome pairing, and used
# "<name>@<surname>.it"-style addresses. This commit replaces them with the
# correctly paired records and normalised "<name>@gmail.com" addresses
# (two of which are deliberately malformed - "luigigmail.com" and
# "federica@gmailcom" - to keep exercising the import-validation-error
# scenario the file is used for), turns the last email (E19) into a real
# mailto: hyperlink, and leaves a stray "Hyperlink" style on E2 as well
# (no functioning link there). It also updates the Dipartimento/Posizione
# (H/I) counters and the current sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A="Giovanni";   B="Moretti";  C="M"; E="giovanni@gmail.com";    H=1;  I=1}
    @{Row=3;  A="Laura";      B="Ricci2";   C="F"; E="laura@gmail.com";       H=2;  I=1}
    @{Row=4;  A="Alessandro"; B="Galli";    C="M"; E="alessandro@gmail.com";  H=6;  I=1}
    @{Row=5;  A="Simona";     B="Bianco";   C="F"; E="simona@gmail.com";      H=1;  I=4}
    @{Row=6;  A="Daniele";    B="Ferrari";  C="M"; E="daniele@gmail.com";     H=5;  I=1}
    @{Row=7;  A="Elena";      B="Rossini";  C="F"; E="elena@gmail.com";       H=1;  I=3}
    @{Row=8;  A="Roberto";    B="Greco";    C="M"; E="roberto@gmail.com";     H=3;  I=1}
    @{Row=9;  A="Silvia";     B="Conti";    C="F"; E="silvia@gmail.com";      H=8;  I=11}
    @{Row=10; A="Marco";      B="Marini";   C="M"; E="marco@gmail.com";       H=1;  I=1}
    @{Row=11; A="Valentina1"; B="Santoro";  C="F"; E="valentina@gmail.com";   H=2;  I=7}
    @{Row=12; A="Nicola";     B="Gallo";    C="M"; E="nicola@gmail.com";      H=1;  I=1}
    @{Row=13; A="Giorgia";    B="Barbieri"; C="F"; E="giorgia@gmail.com";     H=4;  I=1}
    @{Row=14; A="Luigi";      B="Fontana";  C="M"; E="luigigmail.com";        H=11; I=6}
    @{Row=15; A="Elisa";      B="Morelli";  C="F"; E="elisa@gmail.com";       H=1;  I=1}
    @{Row=16; A="Francesco";  B="Riva";     C="M"; E="francesco@gmail.com";   H=1;  I=5}
    @{Row=17; A="Martina";    B="Lombardi"; C="B"; E="martina@gmail.com";     H=1;  I=1}
    @{Row=18; A="Paolo";      B="Colombo";  C="M"; E="paolo@gmail.com";       H=3;  I=9}
    @{Row=19; A="Federica";   B="Caruso";   C="F"; E="federica@gmailcom";     H=9;  I=1}
    @{Row=20; A="Riccardo";   B="Gatti";    C="A"; E="riccardo@gmail.com";    H=1;  I=7}
)

foreach ($row in $data) {
    $ws.Range("A" + $row.Row).Value = $row.A
    $ws.Range("B" + $row.Row).Value = $row.B
    $ws.Range("C" + $row.Row).Value = $row.C
    $ws.Range("E" + $row.Row).Value = $row.E
    $ws.Range("H" + $row.Row).Value = $row.H
    $ws.Range("I" + $row.Row).Value = $row.I
}

# Turn the Federica row's email into a real hyperlink (mailto:), keeping the
# already-set cell text ("federica@gmailcom") as the displayed text.
$ws.Hyperlinks.Add($ws.Range("E19"), "mailto:federica@gmail.com")

# E2 also picked up the "Hyperlink" cell style (underline + theme color)
# without an actual hyperlink being attached to it.
$ws.Range("E2").Style = "Hyperlink"

# Update the active selection shown when the sheet is reopened.
$ws.Range("K12").Select()
